$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New tracker entries for 2025-09-10 (serial date 45910), appended as rows 17-21
$newRows = @(
    @{ A = "G2"; B = "Workout"; C = 45910; D = 0.9705901479276444; E = 0; F = -0.01 },
    @{ A = "G3"; B = "Eat Healthy"; C = 45910; D = 0.9705901479276444; E = 0; F = -0.01 },
    @{ A = "G4"; B = "Read Book"; C = 45910; D = 0.9705901479276444; E = 0; F = -0.01 },
    @{ A = "G5"; B = "Investment Plan"; C = 45910; D = 0.9705901479276444; E = 0; F = -0.01 },
    @{ A = "G6"; B = "Spend 10 Hours without phone"; C = 45910; D = 0.9705901479276444; E = 0; F = -0.01 }
)

$startRow = 17
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B

    $cDate = $ws.Cells.Item($r, 3)
    $cDate.Value = $row.C
    $cDate.NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
}
